$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Battery_Data")
$ws1.Range("B2").Value = 0.53158528266100002
$ws1.Range("B3").Value = 0.39337310916913998
$ws1.Range("B4").Value = 0.0078674621833828007
$ws1.Range("B5").Value = 0.113355856262

$ws2 = $wb.Worksheets.Item("Yearly BRC")
$ws2.Range("B2").Value = 0.05668131650166288
$ws2.Range("B3").Value = 0.056674539760673173
